# Slide 12 ("Optimization: Strength Reduction") - minor modifications:
#  - give the body placeholder (Rectangle 3) an explicit position/size
#  - tweak trailing/leading whitespace around several Consolas code runs

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(4)          # "Rectangle 3" - body placeholder (idx=1)

# --- explicit xfrm on the placeholder (previously inherited, spPr was empty) ---
# Target EMU values: off x=458787 y=1363663, ext cx=8321040 cy=4935537
# Shape.Left/Top/Width/Height are expressed in points (EMU / 12700); PowerPoint
# stores them internally as single-precision floats, so a literal EMU/12700
# for "y" truncates back to 1363662 EMU after the float round-trip. Using a
# point value whose float32 representation still maps back to 1363663 keeps
# the result pixel-exact.
$sh.Left   = 458787 / 12700
$sh.Top    = 107.37504197007874
$sh.Width  = 8321040 / 12700
$sh.Height = 4935537 / 12700

$tf = $sh.TextFrame
$tr = $tf.TextRange

# Several of the edits below live inside a run of text that PowerPoint's COM
# layer *displays* merged together with neighboring same-formatted runs (e.g.
# " ", "inc", " i " all show up as a single ".Runs()" entry). Writing through
# that merged Run object only rewrites the first underlying XML run, so we
# instead address the exact characters with TextRange.Characters(start,len)
# (1-based, absolute offsets into the shape's TextRange) which maps onto the
# real underlying <a:r> run. Apply edits from the highest start offset down
# to the lowest so earlier offsets stay valid as text length changes.

$tr.Characters(344, 30).Text = "(usually smaller and faster)"
$tr.Characters(336, 8).Text  = " XOR EAX "
$tr.Characters(278, 24).Text = "(replace division by 2"
$tr.Characters(271, 7).Text  = " x >> 3  "
$tr.Characters(266, 4).Text  = "x/8        "
$tr.Characters(220, 45).Text = "(replace multiplication by 2 with addition)"
$tr.Characters(214, 6).Text  = " i + i   "
$tr.Characters(171, 3).Text  = " i   "
$tr.Characters(156, 10).Text = "i = i + 1  "
